# NSMB - we have blue shell!!!!
# Apply the edits to the "V4" sheet (sheet1 / first tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# --- Row 1: new header cell K1 ---
$ws.Range("K1").Value = "No bus"

# --- Row 2: new K2 value (bus-level-timing-disabled run) ---
$ws.Range("K2").Value = 2806

# --- Row 3 ---
$ws.Range("K3").Value = 17

# --- Row 4: running total formula ---
$ws.Range("K4").Formula = "=SUM(K2:K3)"

# --- Row 5 ---
$ws.Range("K5").Value = 4146

# --- Row 8: empty styled placeholder cell ---
$ws.Range("K8").Value = ""

# --- Row 9 ---
$ws.Range("B9").Value = 2806
$ws.Range("K9").Value = ""

# --- Row 10 ---
$ws.Range("B10").Value = 3093
$ws.Range("K10").Value = 22200
$ws.Range("L10").Value = "none"

# --- Row 11 ---
$ws.Range("B11").Value = 3151
$ws.Range("K11").Value = 22300
$ws.Range("L11").Value = "none"

# --- Row 12 ---
$ws.Range("B12").Value = 3712
$ws.Range("K12").Value = 22600
$ws.Range("L12").Value = "none"

# --- Row 13 ---
$ws.Range("B13").Value = 4146
$ws.Range("K13").Value = 22900
$ws.Range("L13").Value = "none"

# --- Row 14 ---
$ws.Range("B14").Value = 4664
$ws.Range("K14").Value = 23100
$ws.Range("L14").Value = "none"

# --- Row 15 ---
$ws.Range("K15").Value = 22800
$ws.Range("L15").Value = "none"

# --- Row 16 ---
$ws.Range("K16").Value = 23000

# --- Row 17 ---
$ws.Range("B17").Value = 5057
$ws.Range("K17").Value = 23600
$ws.Range("L17").Value = "none"

# --- Update the active selection to match the saved view ---
$ws.Range("B10").Select()
